# Adding export to docx command in Gui
#
# This reproduces, via the Word object model, the structural changes that
# Word itself performs when the document is re-saved / round-tripped by
# the new "export to docx" path:
#   1. The (only) body paragraph gets an explicit LTR reading order
#      (<w:pPr><w:bidi w:val="0"/></w:pPr>).
#   2. Word's footnote/endnote story parts (word/footnotes.xml,
#      word/endnotes.xml) get minted with their standard separator /
#      continuation-separator boilerplate (plus the matching
#      relationships + content-type overrides), as happens the first
#      time Word touches the footnote/endnote machinery on a document
#      that never had any.
#   3. A couple of custom table styles get tidied up: "connection" picks
#      up a larger default font size and the redundant, never-used
#      "DATAFRAME" table style (a stray duplicate of "Plain Table 1") is
#      removed.

$d = $word.ActiveDocument

# --- 1. Force explicit left-to-right reading order on the lone paragraph ---
$d.Paragraphs(1).Format.ReadingOrder = 0

# --- 2. Mint word/footnotes.xml + word/endnotes.xml -----------------------
# Adding a footnote is the only documented way to get Word to create the
# footnotes/endnotes parts (with their separator/continuationSeparator
# boilerplate); immediately deleting the reference removes the footnote
# itself again while leaving the now-created parts (and their
# relationships / content-type overrides) in place, exactly as an
# add-then-discard edit would.
$r = $d.Paragraphs(1).Range
$r.Collapse(1)
$fn = $d.Footnotes.Add($r, "", "x")
$fn.Reference.Delete()

# --- 3. Table style cleanup -------------------------------------------------
# "connection" style: bump the run font size from 8pt to 12pt.
$connection = $d.Styles("connection")
$connection.Font.Size = 12

# Remove the redundant "DATAFRAME" custom table style (duplicate of
# "Plain Table 1").
$dataframe = $d.Styles("DATAFRAME")
$dataframe.Delete()
